$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 35
$ws.Cells.Item($row, 1).Value = "BonusPower"
$ws.Cells.Item($row, 2).Value = 0.01
$ws.Cells.Item($row, 3).Value = 10
$ws.Cells.Item($row, 4).Value = 6
$ws.Cells.Item($row, 5).Value = 100
$ws.Cells.Item($row, 6).Value = 1
$ws.Cells.Item($row, 7).Value = "win"
